$wb = $excel.ActiveWorkbook

# "Contacts Final Changes - 6th July 2023":
# Update the default/standard user name on the Users sheet from
# "Nicole Bicho" to "Drew Koecher".
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Drew Koecher"

# Leave the workbook with the Users sheet active/selected, matching the
# cursor position the author ended the editing session at.
$ws.Activate() | Out-Null
$ws.Range("B11").Select() | Out-Null
